$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 38 (shifts rows 38:86 down to 39:87)
$ws.Rows.Item(38).Insert()

$ws.Range("A38").Value = 8545
$ws.Range("B38").Value = "한진해모로"

# Update the selection to match the saved view
$ws.Range("B25").Select()
